# Apply corrected results after fixing error estimation and the number of
# projected years (trends not yet rerun after filtering bug fix).

$wb = $excel.ActiveWorkbook

# --- "Trends Status" sheet ---
$ws1 = $wb.Worksheets.Item("Trends Status")
$ws1.Range("C2").Value = 5
$ws1.Range("E2").Value = 11.9
$ws1.Range("C3").Value = 10
$ws1.Range("E3").Value = 23.8
$ws1.Range("C4").Value = 23
$ws1.Range("E4").Value = 54.8
$ws1.Range("E5").Value = 4.8
$ws1.Range("C6").Value = 2
$ws1.Range("E6").Value = 4.8
$ws1.Range("C7").Value = 142

# --- "Species qualification" sheet ---
$ws4 = $wb.Worksheets.Item("Species qualification")
$ws4.Range("C4").Value = 42

# --- "Interannual update - High Pri" sheet ---
$ws5 = $wb.Worksheets.Item("Interannual update - High Pri")
$ws5.Range("B2").Value = 58
$ws5.Range("C2").Value = 56.3
$ws5.Range("D2").Value = 58
$ws5.Range("E2").Value = 73.40000000000001
$ws5.Range("B3").Value = 45
$ws5.Range("C3").Value = 43.7
$ws5.Range("D3").Value = 21
$ws5.Range("E3").Value = 26.6
